$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for numeric-looking Price values so Excel does not
# auto-convert them to numbers (values like "5.34" or "600.90" would
# otherwise be parsed as floating point numbers).
$textCells = @("D5", "D6", "D9", "D11", "D12", "D13", "D14", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D27", "D28", "D29", "D31", "D33", "D34", "D36", "D37", "D39", "D40", "D41", "D43", "D46", "D48", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 31 / Row 33 swap: FirstDigitalUSD <-> EthereumClassic, each with
# refreshed Price / Volume(1h) values.

$ws.Range("D2").Value = "63.730.71"
$ws.Range("E2").Value = "  -1.12%  "
$ws.Range("D3").Value = "3.120.49"
$ws.Range("E3").Value = "  -1.17%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "600.90"
$ws.Range("E5").Value = "  -1.88%  "
$ws.Range("D6").Value = "142.21"
$ws.Range("E6").Value = "  -4.04%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.120.72"
$ws.Range("E8").Value = "  -1.05%  "
$ws.Range("D9").Value = "0.523"
$ws.Range("E9").Value = "  -0.62%  "
$ws.Range("E10").Value = "  -2.37%  "
$ws.Range("D11").Value = "5.34"
$ws.Range("E11").Value = "  -2.34%  "
$ws.Range("D12").Value = "0.465"
$ws.Range("E12").Value = "  -1.70%  "
$ws.Range("D13").Value = "0.0000253"
$ws.Range("E13").Value = "  -2.39%  "
$ws.Range("D14").Value = "34.97"
$ws.Range("E14").Value = "  -1.79%  "
$ws.Range("D15").Value = "3.637.57"
$ws.Range("E15").Value = "  -1.03%  "
$ws.Range("E16").Value = "  +2.52%  "
$ws.Range("D17").Value = "63.757.02"
$ws.Range("E17").Value = "  -0.97%  "
$ws.Range("D18").Value = "3.128.50"
$ws.Range("E18").Value = "  -0.94%  "
$ws.Range("D19").Value = "6.81"
$ws.Range("E19").Value = "  -1.51%  "
$ws.Range("D20").Value = "482.73"
$ws.Range("E20").Value = "  +0.76%  "
$ws.Range("D21").Value = "14.66"
$ws.Range("E21").Value = "  -0.26%  "
$ws.Range("D22").Value = "0.705"
$ws.Range("E22").Value = "  -1.50%  "
$ws.Range("D23").Value = "7.60"
$ws.Range("E23").Value = "  -5.59%  "
$ws.Range("D24").Value = "86.64"
$ws.Range("E24").Value = "  +3.20%  "
$ws.Range("D25").Value = "13.37"
$ws.Range("E25").Value = "  -2.53%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").Value = "2.74"
$ws.Range("E27").Value = "  -4.05%  "
$ws.Range("D28").Value = "8.22"
$ws.Range("E28").Value = "  -4.46%  "
$ws.Range("D29").Value = "6.97"
$ws.Range("E29").Value = "  -2.00%  "
$ws.Range("E30").Value = "  -2.53%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "27.02"
$ws.Range("E31").Value = "  +2.49%  "
$ws.Range("E32").Value = "  -6.92%  "
$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  -0.12%  "
$ws.Range("D34").Value = "2.65"
$ws.Range("E34").Value = "  -2.41%  "
$ws.Range("E35").Value = "  -3.37%  "
$ws.Range("D36").Value = "5.97"
$ws.Range("E36").Value = "  -0.92%  "
$ws.Range("D37").Value = "52.46"
$ws.Range("E37").Value = "  -1.36%  "
$ws.Range("D38").Value = "0.0₃0741"
$ws.Range("E38").Value = "  -7.37%  "
$ws.Range("D39").Value = "2.93"
$ws.Range("E39").Value = "  -10.25%  "
$ws.Range("D40").Value = "436.30"
$ws.Range("E40").Value = "  -5.90%  "
$ws.Range("D41").Value = "0.0394"
$ws.Range("E41").Value = "  -1.84%  "
$ws.Range("E42").Value = "  -0.80%  "
$ws.Range("D43").Value = "8.24"
$ws.Range("E43").Value = "  -1.72%  "
$ws.Range("D44").Value = "2.867.40"
$ws.Range("E44").Value = "  +0.49%  "
$ws.Range("E45").Value = "  -3.92%  "
$ws.Range("D46").Value = "2.20"
$ws.Range("E46").Value = "  -5.84%  "
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("D48").Value = "2.35"
$ws.Range("E48").Value = "  -4.89%  "
$ws.Range("D49").Value = "25.78"
$ws.Range("E49").Value = "  -3.20%  "
$ws.Range("E50").Value = "  -0.62%  "
$ws.Range("D51").Value = "121.58"
$ws.Range("E51").Value = "  +1.19%  "
